$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.081.37"
Set-TextValue $ws.Range("E2") "  -1.73%  "

Set-TextValue $ws.Range("D3") "1.892.25"
Set-TextValue $ws.Range("E3") "  -1.15%  "

Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.14%  "

Set-TextValue $ws.Range("D5") "313.99"
Set-TextValue $ws.Range("E5") "  -0.41%  "

Set-TextValue $ws.Range("E6") "  +0.09%  "

Set-TextValue $ws.Range("D7") "0.5046"
Set-TextValue $ws.Range("E7") "  -0.51%  "

Set-TextValue $ws.Range("D8") "0.3904"
Set-TextValue $ws.Range("E8") "  -1.30%  "

Set-TextValue $ws.Range("D9") "0.09200"
Set-TextValue $ws.Range("E9") "  -6.51%  "

Set-TextValue $ws.Range("D10") "1.126"
Set-TextValue $ws.Range("E10") "  -2.94%  "

Set-TextValue $ws.Range("D11") "41.82"
Set-TextValue $ws.Range("E11") "  -1.18%  "

Set-TextValue $ws.Range("D12") "6.377"
Set-TextValue $ws.Range("E12") "  -2.60%  "

Set-TextValue $ws.Range("D13") "20.79"
Set-TextValue $ws.Range("E13") "  -2.07%  "

Set-TextValue $ws.Range("D14") "1.902.68"
Set-TextValue $ws.Range("E14") "  -0.73%  "

Set-TextValue $ws.Range("D15") "7.279"
Set-TextValue $ws.Range("E15") "  -4.01%  "

Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  +0.17%  "

Set-TextValue $ws.Range("D17") "92.42"
Set-TextValue $ws.Range("E17") "  -1.75%  "

Set-TextValue $ws.Range("D18") "0.00001107"
Set-TextValue $ws.Range("E18") "  -2.99%  "

Set-TextValue $ws.Range("D19") "0.06665"
Set-TextValue $ws.Range("E19") "  +0.01%  "

Set-TextValue $ws.Range("D20") "17.82"
Set-TextValue $ws.Range("E20") "  -1.83%  "

Set-TextValue $ws.Range("D21") "1.001"
Set-TextValue $ws.Range("E21") "  +0.08%  "

Set-TextValue $ws.Range("E22") "  -1.93%  "

Set-TextValue $ws.Range("D23") "28.141.28"

Set-TextValue $ws.Range("D24") "11.39"
Set-TextValue $ws.Range("E24") "  -0.49%  "

Set-TextValue $ws.Range("E25") "  +1.57%  "

Set-TextValue $ws.Range("D26") "2.113.17"
Set-TextValue $ws.Range("E26") "  -1.11%  "

Set-TextValue $ws.Range("E27") "  -7.06%  "

Set-TextValue $ws.Range("B28") "Monero"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "158.22"
Set-TextValue $ws.Range("E28") "  -0.89%  "

Set-TextValue $ws.Range("B29") "EthereumClassic"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D29") "20.84"
Set-TextValue $ws.Range("E29") "  -2.23%  "

Set-TextValue $ws.Range("D30") "126.97"
Set-TextValue $ws.Range("E30") "  -1.46%  "

Set-TextValue $ws.Range("D31") "1.073"
Set-TextValue $ws.Range("E31") "  -2.91%  "

Set-TextValue $ws.Range("E32") "  -1.78%  "

Set-TextValue $ws.Range("D33") "5.603"
Set-TextValue $ws.Range("E33") "  -2.64%  "

Set-TextValue $ws.Range("D34") "3.616"
Set-TextValue $ws.Range("E34") "  -0.85%  "

Set-TextValue $ws.Range("D35") "9.579"
Set-TextValue $ws.Range("E35") "  -3.03%  "

Set-TextValue $ws.Range("D36") "1.349"
Set-TextValue $ws.Range("E36") "  +13.26%  "

Set-TextValue $ws.Range("D37") "0.06604"
Set-TextValue $ws.Range("E37") "  -3.01%  "

Set-TextValue $ws.Range("D38") "0.02400"
Set-TextValue $ws.Range("E38") "  -2.00%  "

Set-TextValue $ws.Range("D39") "0.2206"
Set-TextValue $ws.Range("E39") "  -1.26%  "

Set-TextValue $ws.Range("D40") "1.217"
Set-TextValue $ws.Range("E40") "  -4.47%  "

Set-TextValue $ws.Range("D41") "0.6453"
Set-TextValue $ws.Range("E41") "  +0.05%  "

Set-TextValue $ws.Range("D42") "11.41"
Set-TextValue $ws.Range("E42") "  -3.96%  "

Set-TextValue $ws.Range("D43") "4.959"
Set-TextValue $ws.Range("E43") "  -2.88%  "

Set-TextValue $ws.Range("D44") "1.001"
Set-TextValue $ws.Range("E44") "  +0.09%  "

Set-TextValue $ws.Range("B45") "EnergySwap"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "13.38"
Set-TextValue $ws.Range("E45") "  -2.71%  "

Set-TextValue $ws.Range("B46") "Decentraland"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.6065"
Set-TextValue $ws.Range("E46") "  -0.55%  "

Set-TextValue $ws.Range("D47") "1.303"
Set-TextValue $ws.Range("E47") "  +1.48%  "

Set-TextValue $ws.Range("D48") "3.684"
Set-TextValue $ws.Range("E48") "  -3.49%  "

Set-TextValue $ws.Range("D49") "2.000"
Set-TextValue $ws.Range("E49") "  -2.09%  "

Set-TextValue $ws.Range("D50") "122.17"
Set-TextValue $ws.Range("E50") "  -2.39%  "

Set-TextValue $ws.Range("E51") "  -1.37%  "
